$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Each new value is prefixed with a leading apostrophe so Excel stores it
# as text (matching the original inlineStr cell type) instead of
# auto-converting numeric-looking strings (e.g. "1.00") into numbers.

$ws.Range("D2").Value = "'29.871.62"
$ws.Range("E2").Value = "'  +0.59%  "
$ws.Range("D3").Value = "'1.631.36"
$ws.Range("E3").Value = "'  +1.53%  "
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("D5").Value = "'214.35"
$ws.Range("E5").Value = "'  +0.62%  "
$ws.Range("E6").Value = "'  -0.05%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'28.62"
$ws.Range("E8").Value = "'  +2.00%  "
$ws.Range("E9").Value = "'  +2.06%  "
$ws.Range("E10").Value = "'  +0.80%  "
$ws.Range("E11").Value = "'  +0.50%  "
$ws.Range("D12").Value = "'1.866.61"
$ws.Range("E12").Value = "'  +1.61%  "
$ws.Range("D13").Value = "'1.634.49"
$ws.Range("E13").Value = "'  +1.25%  "
$ws.Range("E14").Value = "'  +3.00%  "
$ws.Range("D15").Value = "'9.28"
$ws.Range("E15").Value = "'  +18.29%  "
$ws.Range("E16").Value = "'  +2.56%  "
$ws.Range("D17").Value = "'29.890.16"
$ws.Range("E17").Value = "'  +0.57%  "
$ws.Range("D18").Value = "'64.18"
$ws.Range("E18").Value = "'  +0.05%  "
$ws.Range("D19").Value = "'242.48"
$ws.Range("E19").Value = "'  +0.42%  "
$ws.Range("D20").Value = "'0.0₃0699"
$ws.Range("E20").Value = "'  +0.16%  "
$ws.Range("E21").Value = "'  +0.06%  "
$ws.Range("D22").Value = "'9.84"
$ws.Range("E22").Value = "'  +4.56%  "
$ws.Range("E23").Value = "'  +2.22%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "'  +0.38%  "
$ws.Range("D25").Value = "'157.56"
$ws.Range("E25").Value = "'  +1.59%  "
$ws.Range("D26").Value = "'15.51"
$ws.Range("E26").Value = "'  +0.24%  "
$ws.Range("E27").Value = "'  +0.95%  "
$ws.Range("D28").Value = "'6.58"
$ws.Range("E28").Value = "'  +1.85%  "
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("E30").Value = "'  +0.98%  "
$ws.Range("E31").Value = "'  +3.98%  "
$ws.Range("E32").Value = "'  +3.77%  "
$ws.Range("E33").Value = "'  -0.42%  "
$ws.Range("D34").Value = "'1.426.66"
$ws.Range("E34").Value = "'  -0.08%  "
$ws.Range("D35").Value = "'1.63"
$ws.Range("E35").Value = "'  +4.49%  "
$ws.Range("E36").Value = "'  +0.03%  "
$ws.Range("D37").Value = "'2.81"
$ws.Range("E37").Value = "'  -2.74%  "
$ws.Range("E38").Value = "'  -0.07%  "
$ws.Range("E39").Value = "'  +0.35%  "
$ws.Range("D40").Value = "'75.47"
$ws.Range("E40").Value = "'  +13.84%  "
$ws.Range("D41").Value = "'0.552"
$ws.Range("E41").Value = "'  +0.60%  "
$ws.Range("E42").Value = "'  +1.92%  "
$ws.Range("E43").Value = "'  +0.90%  "
$ws.Range("E44").Value = "'  -1.57%  "
$ws.Range("D45").Value = "'53.40"
$ws.Range("E45").Value = "'  -6.22%  "
$ws.Range("B46").Value = "'PaxDollar"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "'  +0.02%  "
$ws.Range("B47").Value = "'WEMIXToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.02"
$ws.Range("E47").Value = "'  +3.65%  "
$ws.Range("D48").Value = "'5.36"
$ws.Range("E48").Value = "'  +0.05%  "
$ws.Range("D49").Value = "'1.773.14"
$ws.Range("E49").Value = "'  +1.72%  "
$ws.Range("D50").Value = "'89.06"
$ws.Range("E50").Value = "'  +2.72%  "
$ws.Range("D51").Value = "'0.0₆0110"
$ws.Range("E51").Value = "'  +5.23%  "
